$d = $word.ActiveDocument

# Word (and this Find implementation) walks `wdFindContinue` wraps across
# the *whole* document regardless of which Range invoked `.Find`, so a
# Range scoped to a single paragraph is not a hard search boundary here.
# The reliable approach is: operate on $d.Content (the whole story),
# always match case-sensitively, use wdReplaceOne (1) so only a single
# occurrence is touched per call, and issue the calls in document order
# so each call's *first* remaining match is the intended one.
function ReplaceOnce($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 1) | Out-Null
}

# --- Paragraph 1: Title ---
ReplaceOnce "Read First: Notes for translation" "Citiți mai întâi: Note pentru traducere"

# --- Paragraph 2: "Only translate content..." ---
ReplaceOnce "Only translate content appearing in the columns for " "Traduceți numai conținutul care apare în coloanele pentru"
ReplaceOnce "{Script} " "{Script}"
ReplaceOnce "and " "și"
ReplaceOnce "{On Slide Text}. " " {On Slide Text}. "

# --- Paragraph 3: "Do not translate content in the..." ---
ReplaceOnce "Do not" "Nu"
ReplaceOnce " translate content in the " "traduce conținutul din coloana "
ReplaceOnce "{Animation Notes} " "{Animation Notes}"
ReplaceOnce "column. Leave this column in English. " ". Lasă această coloană în limba engleză. "

# --- Paragraph 4: "Do not translate Lesson headings..." ---
ReplaceOnce "Do not translate Lesson headings" "Nu traduce titlurile lecțiilor"
ReplaceOnce ". Leave them in English. Please make sure the lesson headings remain in English next to “Lesson:”, as this helps the team tracking the video translation. " ". Lasă-le în engleză. Asigură-te că titlurile lecțiilor rămân în limba engleză lângă „Lecție:”, deoarece acest lucru ajută echipa să urmărească traducerea videoclipului. "

# --- Paragraph 5: "Ensure names of tips..." ---
ReplaceOnce "Ensure names of tips are consistent across a module. For example, “Be Calm” should remain in the same form on every slide it appears in. It shouldn’t change to “Stay Calm” or “Being Calm” later on. " "Asigură-te că numele sfaturilor sunt consecvente în cadrul unui modul. De exemplu, „Fii calm” ar trebui să rămână în aceeași formă pe fiecare diapozitiv în care apare. Nu ar trebui să se schimbe în „Rămâi calm” sau „Fii calm” mai târziu. "

# --- Paragraph 6: "Don't include any slashes..." ---
ReplaceOnce "Don’t include any slashes (/) or brackets where they weren’t present in the original script - make sure the script is final and matches the original stylistically." "Nu include nicio bară oblică (/) sau paranteze care nu existau în scenariul original – asigură-te că scenariul este final și că se potrivește stilistic cu originalul."

# --- Paragraph 7: "Do not translate numbers..." ---
ReplaceOnce "Do not translate" "Nu traduce"
ReplaceOnce " numbers appearing in brackets, such as [1] or the word [pause]. " "Nu traduce numerele care apar între paranteze, cum ar fi [1] sau cuvântul [pause]. "

# --- "Changes in mood or behaviour" table cell ---
ReplaceOnce "Changes in mood or behaviour " "Schimbări de dispoziție sau comportament "
ReplaceOnce "Notice secretive  " "Observă dacă este secretos  "
ReplaceOnce "Check in and show you care" "Interesează-te și arătă că îți pasă"

# --- "Also look out for if your child is:" table cell ---
ReplaceOnce "Also look out for if your child is:" "De asemenea, fii atenți dacă copilul:"
ReplaceOnce "Going to places that are unusual for them" "Merge în locuri neobișnuite pentru el"
ReplaceOnce "Getting gifts or money they can’t explain" "🔵 Primește cadouri sau bani despre care nu poate da explicații"
ReplaceOnce "Hiding who they talk to or where they go" "🔵Ascunde persoanele cu care vorbește sau locurile în care merge"
ReplaceOnce "Getting nervous, upset, or secretive when using a phone or going online" "🔵Devine nervos, supărat sau secretos când folosește telefonul sau navighează pe internet"
ReplaceOnce "These signs don’t always mean something is wrong but they are good reasons to check in with your children and show them you care" "Aceste semne nu arată întotdeauna că ceva este în neregulă, dar sunt motive întemeiate pentru a te interesa de starea copiilor și a le arăta că îți pasă"

# --- "Look out:" table cell ---
ReplaceOnce "Look out: " "Fii atent:"
ReplaceOnce "Going to places that are unusual for them" "Merge în locuri neobișnuite pentru el"
ReplaceOnce "Getting gifts or money they can’t explain" "Primește cadouri sau bani despre care nu poate da explicații"
ReplaceOnce "Hiding who they talk to or where they go" "Ascunde persoanele cu care vorbește sau locurile în care merge"
ReplaceOnce "Getting nervous, upset, or secretive when using a phone or going online" "Devine nervos, supărat sau secretos când folosește telefonul sau navighează pe internet"

# --- "Care for Your Children in a Shelter" section ---
ReplaceOnce "Today’s lesson is caring for your children in a shelter." "Lecția de azi este despre îngrijirea copiilor într-un adăpost."
ReplaceOnce "Here are three tips that will help you to support your children if you ever find yourself in a shelter:" "Iată trei sfaturi care te vor ajuta să-ți susții copiii dacă te vei afla vreodată într-un adăpost:"
ReplaceOnce "Care for Your Children in a Shelter" "Îngrijirea copiilor într-un adăpost"

ReplaceOnce "The first tip is to make it familiar." "Primul sfat este să îl faceți familiar."
ReplaceOnce "If possible, take a comfort item from home with them to a shelter" "Dacă este posibil, ia de acasă un lucru care aduce confort"
ReplaceOnce "Decorate the shelter together, with notes or pictures for the walls if you can." "Dacă aveți posibilitatea, decorați împreună adăpostul cu mesaje sau fotografii."
ReplaceOnce "Make it familiar " "Faceți locul familiar "

ReplaceOnce "The second tip is to find the positive." "Al doilea sfat este să găsești aspectele pozitive."
ReplaceOnce "Find something to smile about together. " "Găsiți împreună un motiv pentru a zâmbi. "
ReplaceOnce "Share one good thing that happened at the end of each day with your child - no matter how small. " "La sfârșitul fiecărei zile, împărtășește cu copilul tău un lucru bun care s-a întâmplat, indiferent cât de mic. "
ReplaceOnce "Find the positive " "Găsește aspectele pozitive "

ReplaceOnce "The third tip is to reassure." "Al treilea sfat este acela de a liniști."
ReplaceOnce "Try to reassure your children as often as you can. " "Încearcă să îți liniștești copilul cât poți de des. "
ReplaceOnce "For younger children, you may want to give them a hug or hold them when they are feeling scared. " "Pe copiii mai mici, îi puteți îmbrățișa sau ține în brațe atunci când se simt speriați. "
ReplaceOnce "For older children, let them know that you are here for them and give them the space to talk to you." "În cazul copiilor mai mari, comunică-le că ești alături de ei și oferă-le spațiul necesar ca ei să poată discuta cu tine."
ReplaceOnce "Reassure " "Liniștește "
